$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "interface-id" (row 18) and "management" (row 19) rows entirely;
# this shifts every row below up by two.
$ws.Rows("18:19").Delete()

# "int-virtual-link-desc" (D16) no longer carries the "x" flag next to it.
$ws.Range("E16").ClearContents()

# Match D16's highlight style to the "normal" (non-red) look used elsewhere
# (e.g. A2), by copying formatting only.
$ws.Range("A2").Copy()
$ws.Range("D16").PasteSpecial(-4122)

# Drop the stray notes that used to sit next to "int-virtual-link-desc[]" and
# "ext-cpd []" (now at rows 50 / 51 after the earlier row deletion).
$ws.Range("B50").ClearContents()
$ws.Range("B51").ClearContents()

# "ext-cpd []" (A51) picks up the same normal-look style as D16/A2.
$ws.Range("A2").Copy()
$ws.Range("A51").PasteSpecial(-4122)

$excel.CutCopyMode = 0
